$wb = $excel.ActiveWorkbook
$dp = $wb.Worksheets.Item("DP")

# Add the new "Coin change" worksheet right after the "DP" sheet
$ws = $wb.Worksheets.Add([Type]::Missing, $dp)
$ws.Name = "Coin change"

$ws.Range("A2").Value = "i"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 9
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = 11

$ws.Range("A3").Value = "dp init"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 12
$ws.Range("F3").Value = 12
$ws.Range("G3").Value = 12
$ws.Range("H3").Value = 12
$ws.Range("I3").Value = 12
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = 12
$ws.Range("L3").Value = 12
$ws.Range("M3").Value = 12

$ws.Range("A4").Value = "dp formula"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = "min(12,1+dp[1-1])"
$ws.Range("D4").Value = "min(12,1+dp[2-1],1+dp[2-2])"

$ws.Range("A5").Value = "dp final"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1

$ws.Range("A9").Value = "coins"
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 5

$ws.Columns.Item(1).ColumnWidth = 10.7387387387387
$ws.Columns.Item(3).ColumnWidth = 18.045045045045
$ws.Columns.Item(4).ColumnWidth = 26.7837837837838

$ws.Range("D5").Select()

$ws.Activate()
$excel.ActiveWindow.Zoom = 200
